$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the selection shown in the worksheet view ---
$ws.Range("R22").Select()

# --- New data for rows 17-31, columns J:P ---
# Row 17: header row (bold, matches A17:G17)
$ws.Range("J17").Value = "Task No"
$ws.Range("K17").Value = "Task Duration (mins)"
$ws.Range("L17").Value = "Start Date/Time"
$ws.Range("M17").Value = "End Date/Time"
$ws.Range("N17").Value = "Activity Description"
$ws.Range("O17").Value = "Activity Comments"
$ws.Range("P17").Value = "Task Dependencies"
$ws.Range("J17:P17").Font.Bold = $true

$rows = @(
    @{ Row = 18; TaskNo = 1;  Duration = 10; Desc = "Verify Data Reconciliation Databricks job is running successfully"; Comments = "Job triggered without errors"; Dep = 210 },
    @{ Row = 19; TaskNo = 2;  Duration = 10; Desc = "Validate Data Reconciliation output tables are populated"; Comments = "Record counts available"; Dep = 210 },
    @{ Row = 20; TaskNo = 3;  Duration = 10; Desc = "Verify Quote Cache Deletion pipeline execution"; Comments = "Scheduled job started"; Dep = 220 },
    @{ Row = 21; TaskNo = 4;  Duration = 10; Desc = "Validate Quote Cache deletion for sample quotes"; Comments = "Cache cleared as expected"; Dep = 220 },
    @{ Row = 22; TaskNo = 5;  Duration = 10; Desc = "Validate MQS Event Hub ingestion to Bronze tables"; Comments = "Events landing in Bronze"; Dep = 300 },
    @{ Row = 23; TaskNo = 6;  Duration = 10; Desc = "Validate Silver tables population from Bronze"; Comments = "Silver data available"; Dep = 330 },
    @{ Row = 24; TaskNo = 7;  Duration = 10; Desc = "Validate Gold entities and views"; Comments = "Gold data accessible"; Dep = 360 },
    @{ Row = 25; TaskNo = 8;  Duration = 10; Desc = "Validate CoS ODS & Engineering tables"; Comments = "Tables accessible and refreshed"; Dep = 235 },
    @{ Row = 26; TaskNo = 9;  Duration = 10; Desc = "Validate CoS Analyst views"; Comments = "Analyst views accessible"; Dep = 240 },
    @{ Row = 27; TaskNo = 10; Duration = 10; Desc = "Validate CoS Pricing views"; Comments = "Pricing views accessible"; Dep = 250 },
    @{ Row = 28; TaskNo = 11; Duration = 10; Desc = "Validate Vehicle & Area LRT service deployment"; Comments = "Service up and accessible"; Dep = 430 },
    @{ Row = 29; TaskNo = 12; Duration = 10; Desc = "Validate LRT data load using signed-off file"; Comments = "Data loaded successfully"; Dep = 480 },
    @{ Row = 30; TaskNo = 13; Duration = 5;  Desc = "Verify PROD Databricks PV group permissions"; Comments = "Access validated"; Dep = 460 },
    @{ Row = 31; TaskNo = 14; Duration = 10; Desc = "Monitor PROD jobs/logs for errors post-release"; Comments = "No critical alerts"; Dep = 410 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 10).Value = $r.TaskNo
    $ws.Cells.Item($row, 11).Value = $r.Duration
    $ws.Cells.Item($row, 12).Value = "Post deployment"
    $ws.Cells.Item($row, 13).Value = "Post deployment"
    $ws.Cells.Item($row, 14).Value = $r.Desc
    $ws.Cells.Item($row, 15).Value = $r.Comments
    $ws.Cells.Item($row, 16).Value = $r.Dep
}

# Column P entries are bold across the whole new table (matches existing style used elsewhere for this column)
$ws.Range("P18:P31").Font.Bold = $true
